$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Source data is textual (prices/percentages stored
# as literal strings), so force Text format first to stop Excel from
# auto-coercing number-looking strings (e.g. "0.610" -> 0.61).
$updates = [ordered]@{
    'D2' = '33.841.85'
    'E2' = '  -0.86%  '
    'D3' = '1.780.91'
    'E3' = '  -1.16%  '
    'E4' = '  +0.09%  '
    'D5' = '224.58'
    'E5' = '  +0.66%  '
    'D6' = '0.547'
    'E6' = '  -1.02%  '
    'E7' = '  +0.08%  '
    'D8' = '31.92'
    'E8' = '  -1.53%  '
    'E9' = '  +1.09%  '
    'E10' = '  -5.47%  '
    'E11' = '  +1.13%  '
    'D12' = '2.036.78'
    'E12' = '  -1.17%  '
    'E13' = '  +4.67%  '
    'D14' = '1.780.79'
    'E14' = '  -1.36%  '
    'D15' = '33.861.42'
    'E15' = '  -0.88%  '
    'D16' = '0.610'
    'E16' = '  -3.58%  '
    'D17' = '4.14'
    'E17' = '  -1.98%  '
    'D18' = '66.71'
    'E18' = '  -2.48%  '
    'D19' = '238.65'
    'E19' = '  -3.66%  '
    'D20' = '0.0₃0774'
    'E20' = '  -1.88%  '
    'E21' = '  +0.03%  '
    'D22' = '10.57'
    'E22' = '  -2.94%  '
    'E23' = '  -2.31%  '
    'E24' = '  -2.09%  '
    'D25' = '160.54'
    'E25' = '  +0.24%  '
    'E26' = '  -0.81%  '
    'D27' = '16.09'
    'E27' = '  -3.04%  '
    'E28' = '  -0.56%  '
    'E29' = '  +0.19%  '
    'D30' = '1.22'
    'E30' = '  +0.94%  '
    'E31' = '  -2.90%  '
    'D32' = '3.60'
    'E32' = '  -3.65%  '
    'E33' = '  +0.08%  '
    'E34' = '  -1.66%  '
    'D35' = '1.389.99'
    'E35' = '  -1.81%  '
    'E36' = '  -1.68%  '
    'E37' = '  -1.79%  '
    'E38' = '  -1.29%  '
    'D39' = '2.26'
    'E39' = '  +4.99%  '
    'D40' = '2.39'
    'E40' = '  +0.66%  '
    'D41' = '78.63'
    'E41' = '  -2.13%  '
    'D42' = '0.912'
    'E42' = '  -3.74%  '
    'D43' = '13.58'
    'E43' = '  +13.57%  '
    'E44' = '  -3.02%  '
    'D45' = '0.0₆0139'
    'E45' = '  +11.80%  '
    'D46' = '0.0507'
    'E46' = '  +2.84%  '
    'E47' = '  +3.19%  '
    'B48' = 'FraxShare'
    'C48' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D48' = '5.86'
    'E48' = '  -1.57%  '
    'B49' = 'Quant'
    'C49' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D49' = '106.75'
    'E49' = '  -0.79%  '
    'D50' = '1.938.31'
    'E50' = '  -1.30%  '
    'E51' = '  +0.13%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
